# Updated 2D training schedules, no break screen
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "praclen" column (I) from 4 to 5 for the existing training trials
$ws.Range("I2:I5").Value = 5

# Append a new training trial as row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 2
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5
$ws.Cells.Item(6, 8).Value = 21
$ws.Cells.Item(6, 9).Value = 5
$ws.Cells.Item(6, 10).Value = "train_dim2_1"

# Move the active selection onto the newly added row, column A
$ws.Range("A6").Select()
